$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = [double]"9.843764488968567E-08"
$ws.Range("E2").Value = [double]"9.843764488968567E-08"

# Row 3
$ws.Range("D3").Value = [double]"0.9999999998453513"
$ws.Range("E3").Value = [double]"0.9999999998453513"

# Row 4
$ws.Range("D4").Value = [double]"0.007991009692874508"
$ws.Range("E4").Value = [double]"0.007991009692874508"

# Row 5
$ws.Range("D5").Value = [double]"0.002293210012805644"
$ws.Range("E5").Value = [double]"0.002293210012805644"

# Row 6
$ws.Range("D6").Value = [double]"1.652004254434257E-05"
$ws.Range("E6").Value = [double]"1.652004254434257E-05"

# Row 7
$ws.Range("D7").Value = [double]"0.9999999999864893"
$ws.Range("E7").Value = [double]"1.351074807587338E-11"

# Row 8
$ws.Range("D8").Value = [double]"2.336460184481626E-09"
$ws.Range("E8").Value = [double]"0.9999999976635399"

# Row 9
$ws.Range("D9").Value = [double]"1.707736348573142E-09"
$ws.Range("E9").Value = [double]"0.9999999982922636"

# Row 10
$ws.Range("D10").Value = [double]"6.342297046439087E-11"
$ws.Range("E10").Value = [double]"0.9999999999365771"

# Row 11
$ws.Range("D11").Value = [double]"0.07513236942895306"
$ws.Range("E11").Value = [double]"0.9248676305710469"
$ws.Range("F11").Value = [double]"8.873262405395508"
